$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# Plain single-dot decimal strings in column D are prefixed with a leading
# apostrophe so Excel stores them as TEXT (matching the sheet's existing
# inline-string convention) instead of auto-converting them to numbers.

$ws.Range('D2').Value = '41.873.22'
$ws.Range('E2').Value = '  +4.29%  '

$ws.Range('D3').Value = '2.275.98'
$ws.Range('E3').Value = '  +2.33%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''305.53'
$ws.Range('E5').Value = '  +4.08%  '

$ws.Range('D6').Value = '''92.66'
$ws.Range('E6').Value = '  +5.49%  '

$ws.Range('E7').Value = '  +3.89%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('E9').Value = '  +3.33%  '

$ws.Range('D10').Value = '''32.82'
$ws.Range('E10').Value = '  +7.09%  '

$ws.Range('D11').Value = '''53.84'
$ws.Range('E11').Value = '  +5.79%  '

$ws.Range('E12').Value = '  +2.25%  '

$ws.Range('E13').Value = '  +1.48%  '

$ws.Range('D14').Value = '''6.65'
$ws.Range('E14').Value = '  +3.41%  '

$ws.Range('D15').Value = '2.628.33'
$ws.Range('E15').Value = '  +2.29%  '

$ws.Range('D16').Value = '''14.33'
$ws.Range('E16').Value = '  +3.58%  '

$ws.Range('D17').Value = '2.298.67'
$ws.Range('E17').Value = '  +3.64%  '

$ws.Range('E18').Value = '  +3.87%  '

$ws.Range('D19').Value = '41.810.32'
$ws.Range('E19').Value = '  +4.28%  '

$ws.Range('D20').Value = '''12.36'
$ws.Range('E20').Value = '  +9.83%  '

$ws.Range('D21').Value = '0.0₃0909'
$ws.Range('E21').Value = '  +2.14%  '

$ws.Range('E22').Value = '  +2.76%  '

$ws.Range('E23').Value = '  +2.26%  '

$ws.Range('D24').Value = '''242.89'
$ws.Range('E24').Value = '  +2.96%  '

$ws.Range('E25').Value = '  +5.39%  '

$ws.Range('E27').Value = '  +5.55%  '

$ws.Range('D28').Value = '''24.33'
$ws.Range('E28').Value = '  +4.94%  '

$ws.Range('E29').Value = '  +3.13%  '

$ws.Range('D30').Value = '''2.08'
$ws.Range('E30').Value = '  +0.73%  '

$ws.Range('D31').Value = '''34.28'
$ws.Range('E31').Value = '  +7.44%  '

$ws.Range('D32').Value = '''158.98'
$ws.Range('E32').Value = '  +0.28%  '

$ws.Range('E33').Value = '  +0.00%  '

$ws.Range('E34').Value = '  +4.52%  '

$ws.Range('D35').Value = '''0.0751'
$ws.Range('E35').Value = '  +4.96%  '

$ws.Range('D36').Value = '''3.04'
$ws.Range('E36').Value = '  +0.93%  '

$ws.Range('D37').Value = '''17.10'
$ws.Range('E37').Value = '  +9.18%  '

$ws.Range('E38').Value = '  +1.95%  '

$ws.Range('E39').Value = '  +2.79%  '

$ws.Range('E40').Value = '  +4.93%  '

$ws.Range('E41').Value = '  +3.29%  '

$ws.Range('E42').Value = '  +4.64%  '

$ws.Range('D43').Value = '2.071.21'
$ws.Range('E43').Value = '  -0.73%  '

$ws.Range('D44').Value = '''19.61'
$ws.Range('E44').Value = '  +1.69%  '

$ws.Range('E45').Value = '  +3.24%  '

$ws.Range('D46').Value = '''10.36'
$ws.Range('E46').Value = '  +2.99%  '

$ws.Range('E47').Value = '  +5.80%  '

$ws.Range('E48').Value = '  +7.48%  '

$ws.Range('D49').Value = '''73.06'
$ws.Range('E49').Value = '  +7.42%  '

$ws.Range('D50').Value = '''1.53'
$ws.Range('E50').Value = '  +3.42%  '

$ws.Range('E51').Value = '  +3.60%  '
